# Update the project schema to reflect the new feature roadmap.
#
# The single existing "Template" sheet becomes the per-month sheet
# ("Month-Template"), and a brand-new blank sheet ("Overall-Template") is
# inserted in front of it to hold the (future) roadmap/overall view. The
# new sheet becomes tab 1, the renamed month sheet becomes tab 2 and stays
# the active/selected sheet, with D3 selected.

$wb = $excel.ActiveWorkbook

# Grab the existing (only) sheet - this is what will become "Month-Template".
$monthSheet = $wb.Worksheets.Item("Template")

# Insert a brand-new blank worksheet immediately before it.
$overallSheet = $wb.Worksheets.Add($monthSheet)
$overallSheet.Name = "Overall-Template"

# Rename the original sheet now that the new one has been placed ahead of it.
$monthSheet = $wb.Worksheets.Item("Template")
$monthSheet.Name = "Month-Template"

# Keep "Month-Template" the active tab, with D3 selected.
$monthSheet.Activate()
[void]$monthSheet.Range("D3").Select()
